$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$lastRow = $used.Rows.Count
$lastCol = $used.Columns.Count

# 1. Rename header row: "<label>_old" -> "<label>_FV2410", "<label>_new" -> "<label>_FV2504"
#    (the new `<formatversion>` naming scheme replaces the generic old/new suffixes)
for ($c = 1; $c -le $lastCol; $c++) {
    $cell = $ws.Cells.Item(1, $c)
    $v = $cell.Value()
    if ($v -ne $null) {
        $newVal = $v -replace '_old$', '_FV2410' -replace '_new$', '_FV2504'
        if ($newVal -ne $v) {
            $cell.Value = $newVal
        }
    }
}

# 2. Freeze the header row so it stays visible while scrolling
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true

# 3. Turn the data range into a proper Excel Table (ListObject) with an AutoFilter
$tbl = $ws.ListObjects.Add(1, $used, $null, 1)
$tbl.Name = "Table1"
$tbl.TableStyle = ""

Write-Host "Renamed headers, froze header row, and created Table1 over" $used.Address()
